$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.557511329650879
$ws.Range("B1").Value = 6.404435157775879
$ws.Range("C1").Value = 8.778757095336914
$ws.Range("D1").Value = 9.232038497924805
$ws.Range("E1").Value = 1.456969141960144
